$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$failMsg = "실패: Step '텍스트 검색' failed: TextExtractor.find_text() got an unexpected keyword argument 'confidence'"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 6).Value = $failMsg
}
